$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial value that was updated
# from 45190 (2023-09-21) to 45192 (2023-09-23) for every data row (2-90).
for ($row = 2; $row -le 90; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45190) {
        $cell.Value2 = 45192
    }
}
